$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 19460.04417134706
$ws.Range("D3").Value = 905.1407206307279

$ws.Range("B4").Value = 8299.418386005162
$ws.Range("D4").Value = 603.2118482439402

$ws.Range("B5").Value = 2748.034000000001

$ws.Range("B6").Value = 13043.009
$ws.Range("D6").Value = 205.001

$ws.Range("B7").Value = 17317.03250000001
$ws.Range("D7").Value = 1165

$ws.Range("B8").Value = 26749.9910000001
$ws.Range("D8").Value = 1285

$ws.Range("B9").Value = 43289.39100000006
$ws.Range("D9").Value = 6614.002

$ws.Range("F10").Value = 10228199216.46202

$ws.Range("G11").Value = 0.8167863866874108

$ws.Range("F12").Value = 475741962.7639999
$ws.Range("G12").Value = 0.04651277831959957

$ws.Range("G13").Value = 0.1367008349929897
